# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" sheet (duplicated from "2022-Q3" so that
#    sheetPr / pageMargins / column styles come along for free) right
#    after the "总计" summary sheet, then overwrite its contents with the
#    2022-Q4 fund-holding rows.
# 2. Insert a new row into "总计" for the 2022-Q4 summary figures, pushing
#    the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item("2022-Q3")

# --- Step 1: build the new "2022-Q4" sheet ---------------------------------
$q3.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# "2022-Q3" only has one data row (row 2); the new sheet needs three, so
# insert two more rows below it, copying row 2's formatting into them.
$q4.Range("A3:H3").Insert(-4121)
$q4.Range("A3:H3").Insert(-4121)
$q4.Range("A2:H2").Copy()
$q4.Range("A3:H3").PasteSpecial(-4122)
$q4.Range("A4:H4").PasteSpecial(-4122)
$q4.Cells.Item(1,1).Copy()
$q4.Application.CutCopyMode = $false

# Header row
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Helper: write a value into a cell while forcing text storage (so numeric
# looking strings like "9.60" or "0.3638" keep their literal text, matching
# the source data which stores these figures as text, not numbers), without
# leaving a residual NumberFormat-driven cell style behind.
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$q4.Cells.Item(2,1).Value = 0
Set-TextCell $q4.Cells.Item(2,2) "016250"
$q4.Cells.Item(2,3).Value = "华夏远见成长一年持有混合A"
Set-TextCell $q4.Cells.Item(2,4) "9.60"
Set-TextCell $q4.Cells.Item(2,5) "88.62"
Set-TextCell $q4.Cells.Item(2,6) "3.79"
Set-TextCell $q4.Cells.Item(2,7) "0.3638"
$q4.Cells.Item(2,8).Value = 7

# Row 3
$q4.Cells.Item(3,1).Value = 1
Set-TextCell $q4.Cells.Item(3,2) "016251"
$q4.Cells.Item(3,3).Value = "华夏远见成长一年持有混合C"
Set-TextCell $q4.Cells.Item(3,4) "2.97"
Set-TextCell $q4.Cells.Item(3,5) "88.62"
Set-TextCell $q4.Cells.Item(3,6) "3.79"
Set-TextCell $q4.Cells.Item(3,7) "0.1126"
$q4.Cells.Item(3,8).Value = 7

# Row 4
$q4.Cells.Item(4,1).Value = 2
Set-TextCell $q4.Cells.Item(4,2) "005444"
$q4.Cells.Item(4,3).Value = "光大保德信多策略精选18个月定期开放灵活配置混合"
Set-TextCell $q4.Cells.Item(4,4) "0.54"
Set-TextCell $q4.Cells.Item(4,5) "29.43"
Set-TextCell $q4.Cells.Item(4,6) "1.94"
Set-TextCell $q4.Cells.Item(4,7) "0.0105"
$q4.Cells.Item(4,8).Value = 7

# --- Step 2: add the 2022-Q4 row to "总计" -----------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)
$total.Application.CutCopyMode = $false

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 3
$total.Cells.Item(2,4).Value = 0.49

# The "A" column is a 0-based row index (pandas-style); renumber rows 3-8
# now that everything shifted down by one.
for ($r = 3; $r -le 8; $r++) {
    $total.Cells.Item($r,1).Value = $r - 2
}
